$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to be treated/stored as text so Excel does not
    # auto-convert numeric-looking strings (e.g. "230.82") into numbers,
    # then restore the default ("Normal") cell style so no stray style
    # reference is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.416.16"
$ws.Range("E2").Value = "  +4.09%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.217.97"
$ws.Range("E3").Value = "  +2.63%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "230.82"
$ws.Range("E5").Value = "  +1.72%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.623"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "60.88"
$ws.Range("E7").Value = "  -3.34%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.402"
$ws.Range("E9").Value = "  +2.64%  "

# Row 10 - OKB
Set-TextValue $ws.Range("D10") "58.78"
$ws.Range("E10").Value = "  +0.55%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0891"
$ws.Range("E11").Value = "  +5.61%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.13%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "2.546.54"
$ws.Range("E13").Value = "  +2.64%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "15.57"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "21.66"
$ws.Range("E15").Value = "  -0.76%  "

# Row 16 - Polygon
Set-TextValue $ws.Range("D16") "0.796"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 - Polkadot
Set-TextValue $ws.Range("D17") "5.54"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.217.40"
$ws.Range("E18").Value = "  +2.64%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "41.266.73"
$ws.Range("E19").Value = "  +3.83%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "72.77"
$ws.Range("E20").Value = "  +1.35%  "

# Row 21 - ShibaInu
Set-TextValue $ws.Range("D21") "0.0₃0893"
$ws.Range("E21").Value = "  +5.58%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "6.04"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "249.70"
$ws.Range("E23").Value = "  +8.63%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.12%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.39"
$ws.Range("E25").Value = "  +0.91%  "

# Row 26 - Toncoin
Set-TextValue $ws.Range("D26") "2.32"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.41"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "168.12"
$ws.Range("E28").Value = "  -2.45%  "

# Row 29 - Kaspa
Set-TextValue $ws.Range("D29") "0.139"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30 - EthereumClassic
Set-TextValue $ws.Range("D30") "19.91"
$ws.Range("E30").Value = "  +0.70%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.13%  "

# Row 32 - WEMIXToken
Set-TextValue $ws.Range("D32") "2.63"
$ws.Range("E32").Value = "  -2.50%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +0.11%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "4.97"
$ws.Range("E34").Value = "  +6.12%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  +0.58%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.0622"
$ws.Range("E36").Value = "  +0.72%  "

# Row 37 - THORChain
$ws.Range("E37").Value = "  -5.63%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "3.68"
$ws.Range("E38").Value = "  -1.85%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -1.29%  "

# Row 40 - TerraClassic
Set-TextValue $ws.Range("D40") "0.000247"
$ws.Range("E40").Value = "  +29.40%  "

# Row 41 - BinanceUSD
Set-TextValue $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - was FTXToken, now VeChain
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0237"
$ws.Range("E42").Value = "  +3.96%  "

# Row 43 - was VeChain, now FTXToken
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D43") "4.81"
$ws.Range("E43").Value = "  -2.60%  "

# Row 44 - FraxShare
Set-TextValue $ws.Range("D44") "8.56"
$ws.Range("E44").Value = "  +8.56%  "

# Row 45 - Cronos
Set-TextValue $ws.Range("D45") "0.0981"
$ws.Range("E45").Value = "  +6.34%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "99.03"
$ws.Range("E46").Value = "  -3.50%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  -0.47%  "

# Row 48 - Maker
Set-TextValue $ws.Range("D48") "1.463.82"
$ws.Range("E48").Value = "  -3.23%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "16.52"
$ws.Range("E49").Value = "  -6.11%  "

# Row 50 - HuobiToken
Set-TextValue $ws.Range("D50") "2.78"
$ws.Range("E50").Value = "  -0.83%  "

# Row 51 - was ARBITRUM, now MultiversX
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "52.31"
$ws.Range("E51").Value = "  +3.98%  "
